# Base class configuration for scottish monthly scenarios:
# unify the TestResultExcelFilePath values on the
# "ProcessPayrollForNIWeekly" sheet so every row points at the 201819
# (current tax year) results workbook instead of the stale 201718 one.

$wb = $excel.ActiveWorkbook

$newPath = "F:\\Automation_TestResults\\Payroll_Tax_NI_Directors_TestReports 201819\\201819 Payroll National Insurance calculation Test result.xlsx"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsProcess.Range("H2").Value = $newPath
$wsProcess.Range("H3").Value = $newPath
$wsProcess.Range("H4").Value = $newPath
$wsProcess.Range("H5").Value = $newPath
$wsProcess.Range("H6").Value = $newPath

# Rows 7-10 already referenced the 201819 report; once all rows share the
# exact same wrapped text, Excel re-applies the same row height (45pt) that
# rows 3-6 already use instead of the old shorter 30pt height.
$wsProcess.Rows.Item(7).RowHeight = 45
$wsProcess.Rows.Item(8).RowHeight = 45
$wsProcess.Rows.Item(9).RowHeight = 45
$wsProcess.Rows.Item(10).RowHeight = 45

# Restore cursor / selection positions on each sheet as left by the author.
$wsNI = $wb.Worksheets.Item("NIWeeklyCat_J")
$wsNI.Select()
$wsNI.Range("B19").Select()

$wsProcess.Select()
$wsProcess.Range("H3:H10").Select()

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Select()
$wsReports.Range("K7").Select()
